$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'34.405.65"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.25%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.787.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -2.08%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.04%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'224.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -2.28%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.552"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -3.65%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.10%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'32.71"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +3.72%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.282"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -2.35%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.0661"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -2.53%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0932"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.07%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'2.043.85"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -2.12%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'11.01"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +6.92%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'1.798.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.23%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.634"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -2.76%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'34.366.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.05%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'4.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.51%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'69.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.91%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'255.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -1.10%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.0₃0746"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.41%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'1.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.24%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'10.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -2.57%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'4.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -3.07%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -4.12%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'157.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.25%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'16.44"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.00%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'7.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -2.04%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.114"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -3.67%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.10%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'3.78"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -3.04%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.0514"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -2.09%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'1.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -2.02%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'3.59"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.23%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.89"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +4.13%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.454.44"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -5.27%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'1.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -2.53%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.631"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.01%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.0189"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -1.13%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'83.35"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.16%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'2.85"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +1.45%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'2.35"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.46%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.891"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -2.61%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'2.07"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -2.46%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.0508"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -3.90%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.06"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -2.15%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'5.87"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.26%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'1.942.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.88%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'12.21"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.22%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.10%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'99.27"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.54%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'50.17"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -3.10%  "
$ws.Range("E51").Style = "Normal"
